$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Created By:" bot name update
$ws.Range("B2").Value = "UiPathBot"

# Remove the extra header/detail columns (F:I) and the column-header row (row 4)
# content -- the template no longer carries the Item/Facility/Buyer/... table
# header, just the plain 5-column (A:E) blank grid that starts at row 4.
$ws.Range("F1:I77").Clear()
$ws.Range("A4:I4").Clear()

# Re-apply the plain body-row formatting (same as row 5) to the now-blank
# header row so it matches the rest of the grid.
$ws.Range("A5:E5").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Rows.Item(4).AutoFit()
$excel.CutCopyMode = 0

# Drop the final spare blank row (77) -- grid now ends at row 76.
$ws.Rows.Item(77).Delete()

# Selection moves to B3 (empty separator row) after the edit.
$ws.Range("B3").Select()
